$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.138.87"
$ws.Range("E2").Value = "  -1.06%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.513.46"
$ws.Range("E3").Value = "  +0.29%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "610.22"
$ws.Range("E5").Value = "  +0.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.21"
$ws.Range("E6").Value = "  -1.44%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.511.11"
$ws.Range("E7").Value = "  +0.26%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.480"
$ws.Range("E9").Value = "  -1.55%  "

$ws.Range("E10").Value = "  -1.02%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "8.05"
$ws.Range("E11").Value = "  +6.95%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.423"
$ws.Range("E12").Value = "  -1.63%  "

$ws.Range("E13").Value = "  +0.56%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.103.91"
$ws.Range("E14").Value = "  +0.44%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.51"
$ws.Range("E15").Value = "  -1.18%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.512.99"
$ws.Range("E16").Value = "  +0.67%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.078.47"
$ws.Range("E17").Value = "  -0.99%  "

$ws.Range("E18").Value = "  -0.17%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.93"
$ws.Range("E19").Value = "  +9.59%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.37"
$ws.Range("E20").Value = "  -1.99%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.41"
$ws.Range("E21").Value = "  +0.41%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "436.72"
$ws.Range("E22").Value = "  -1.88%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.609"
$ws.Range("E23").Value = "  -2.72%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.86"
$ws.Range("E24").Value = "  +1.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.651.66"
$ws.Range("E25").Value = "  +0.49%  "

$ws.Range("E27").Value = "  -3.86%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.82"
$ws.Range("E28").Value = "  -1.54%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.25"
$ws.Range("E29").Value = "  -4.82%  "

$ws.Range("E30").Value = "  +0.57%  "

$ws.Range("E31").Value = "  -4.37%  "

$ws.Range("E32").Value = "  -0.07%  "

$ws.Range("E33").Value = "  -1.90%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.60"
$ws.Range("E34").Value = "  +0.02%  "

$ws.Range("E35").Value = "  -3.13%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.81"
$ws.Range("E36").Value = "  -2.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.05"
$ws.Range("E37").Value = "  +0.66%  "

$ws.Range("E38").Value = "  +0.00%  "

$ws.Range("E39").Value = "  +0.07%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "175.97"
$ws.Range("E40").Value = "  +0.62%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0899"
$ws.Range("E41").Value = "  -0.44%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.42"
$ws.Range("E42").Value = "  -0.32%  "

$ws.Range("E43").Value = "  -9.99%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.897"
$ws.Range("E44").Value = "  -0.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "46.29"
$ws.Range("E45").Value = "  -1.34%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.13"
$ws.Range("E46").Value = "  -8.03%  "

$ws.Range("E47").Value = "  -3.65%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.48"
$ws.Range("E48").Value = "  -1.90%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.46"
$ws.Range("E49").Value = "  -1.64%  "

$ws.Range("E50").Value = "  -0.48%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.248"
$ws.Range("E51").Value = "  -1.79%  "
